$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Give C25:C29 the same date format already used by C2:C24 (xlPasteFormats = -4122)
# before filling in the values below.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C25:C29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Fill in new camera-count summary rows (25-29) ---
# Row 25 already has A25 ("Kenttarova Spruce Ground"); add the rest of the row.
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = 43501
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 1

$ws.Range("A26").Value = "Kenttarova Spruce Ground"
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 43551
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0

$ws.Range("A27").Value = "Kenttarova Spruce Ground"
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = 43592
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0

$ws.Range("A28").Value = "Kenttarova Spruce Ground"
$ws.Range("B28").Value = 4
$ws.Range("C28").Value = 43773
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 1

$ws.Range("A29").Value = "Kenttarova Spruce Ground"
$ws.Range("B29").Value = 5
$ws.Range("C29").Value = 43864

# --- Remove the old list of site names that used to live in column A (rows 33-63) ---
$ws.Range("A33:A63").ClearContents() | Out-Null

# --- Update the saved view state for Sheet2 ---
$ws.Range("C27").Select() | Out-Null
$excel.ActiveWindow.Zoom = 115
